# Applies the "average with safety stocks" edit:
#  - Productdata sheet: divide InventoryCosts (D), BackorderCosts (F) and
#    LostSale (I) columns, rows 2-11, by 2500.
#  - ForcastedStandardDeviation sheet: zero out the standard deviation
#    values (columns B-E) for the last three rows (9-11), since those
#    buckets no longer carry demand uncertainty once safety stocks are
#    averaged in.

$wb = $excel.ActiveWorkbook

$productData = $wb.Worksheets.Item("Productdata")
$divisor = 2500

for ($row = 2; $row -le 11; $row++) {
    foreach ($col in @("D", "F", "I")) {
        $cell = $productData.Range("$col$row")
        $cell.Value2 = $cell.Value2 / $divisor
    }
}

$stdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")

for ($row = 9; $row -le 11; $row++) {
    foreach ($col in @("B", "C", "D", "E")) {
        $stdDev.Range("$col$row").Value = 0
    }
}
